# Fruta / hortaliza, semanal
# Insert two new weekly observations (rows) above the existing row 37 data,
# shifting the rest of the table down by two rows, then populate the two
# new rows with their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 37 - everything from the old row 37
# onward (through row 60) shifts down to rows 39-62.
$ws.Rows("37:38").Insert()

# New row 37
$ws.Cells.Item(37, 1).Value = 6
$ws.Cells.Item(37, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(37, 3).Value = "Metropolitana"
$ws.Cells.Item(37, 4).Value = 44566
$ws.Cells.Item(37, 5).Value = 13
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100101
$ws.Cells.Item(37, 8).Value = "Berries"
$ws.Cells.Item(37, 9).Value = 100101008
$ws.Cells.Item(37, 10).Value = "Mora"
$ws.Cells.Item(37, 11).Value = "Sin especificar"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 250
$ws.Cells.Item(37, 14).Value = 6000
$ws.Cells.Item(37, 15).Value = 6000
$ws.Cells.Item(37, 16).Value = 6000
$ws.Cells.Item(37, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(37, 18).Value = "Provincia de Linares"
$ws.Cells.Item(37, 19).Value = 3000
$ws.Cells.Item(37, 20).Value = 2

# New row 38
$ws.Cells.Item(38, 1).Value = 6
$ws.Cells.Item(38, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(38, 3).Value = "Metropolitana"
$ws.Cells.Item(38, 4).Value = 44566
$ws.Cells.Item(38, 5).Value = 13
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100101
$ws.Cells.Item(38, 8).Value = "Berries"
$ws.Cells.Item(38, 9).Value = 100101008
$ws.Cells.Item(38, 10).Value = "Mora"
$ws.Cells.Item(38, 11).Value = "Sin especificar"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 250
$ws.Cells.Item(38, 14).Value = 5600
$ws.Cells.Item(38, 15).Value = 5600
$ws.Cells.Item(38, 16).Value = 5600
$ws.Cells.Item(38, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(38, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(38, 19).Value = 2800
$ws.Cells.Item(38, 20).Value = 2
